$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
$sysRows = @(3,6,10,11,12,13,14,15,17,18,19,20,21,22,30,33,37,38,39,40,41,42,44,45,46,47,48,49,57,60,64,65,66,67,68,69,71,72,73,74,75,76,86,87,88,89,93,95,96,97,99,102,112,113,114,115,119,121,122,123,125,128,138,139,140,141,145,147,148,149,151,154)
foreach ($r in $sysRows) {
    $ws.Range("G$r").Value = "dnasr281@gmail.com, System"
}

# Rows where "backup@backdoor.com, System" -> "System, backup@backdoor.com"
$backupRows = @(4,31,58)
foreach ($r in $backupRows) {
    $ws.Range("G$r").Value = "System, backup@backdoor.com"
}

Write-Output "Done: updated $($sysRows.Count) System/dnasr281 rows and $($backupRows.Count) backup rows"
